$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.135.53"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.128.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.47"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.42"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.45%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000249"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.61"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.648.19"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.101.28"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.129.63"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.39"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "493.18"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.711"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.92"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +5.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.29"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.34"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.95"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.35"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.76"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0950"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.02%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.91"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.978"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "46.71"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.07"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.15"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.57"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.840.37"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "386.42"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.61"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.61%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "135.83"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.29%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.25%  "
